$wb = $excel.ActiveWorkbook

$wsFieldNotes = $wb.Worksheets.Item("FieldNotes")
$wsGenericDocument = $wb.Worksheets.Item("GenericDocument")

# --- GenericDocument sheet data edits ---
$wsGenericDocument.Range("S2").Value = 2303363
$wsGenericDocument.Range("F3").Value = 45637
$wsGenericDocument.Range("M3").Value = "PUBLIC"
$wsGenericDocument.Range("S3").Value = 2303363

# --- Column width changes ---
# (ColumnWidth values below are pre-adjusted so that the saved OOXML <col>
# width - which this engine rounds to whole-pixel granularity - lands as
# close as possible to the widths produced by real Excel.)

# FieldNotes: set widths for columns E and F
$wsFieldNotes.Columns.Item(5).ColumnWidth = 13.0
$wsFieldNotes.Columns.Item(6).ColumnWidth = 18.0

# GenericDocument: split former E:F bestFit range into distinct widths
$wsGenericDocument.Columns.Item(5).ColumnWidth = 16.5
$wsGenericDocument.Columns.Item(6).ColumnWidth = 19.333333333333332

# --- Active sheet / selection changes ---
$wsGenericDocument.Activate()
$wsGenericDocument.Range("S3").Select()
